# Applies the commit "add synthesized alu with abc and update results":
#  - reformats the P:U (PCF/Synopsys) block on several rows to the
#    comma "#,##0" numeric style already used elsewhere in the sheet
#    (style 7 for the plain number/formula cells, style 6 -- bold -- for
#    the Weight/"U" column), extending this formatting to previously
#    blank rows (7, 10, 11, 14) as well;
#  - adds a new row 13 style placeholder cell (A13) and a new row 14
#    set of placeholder cells matching the surrounding rows;
#  - fills in a brand new "alu" data block (columns D:I) for rows 15
#    and 16, mirroring the formulas used by every other row of the
#    table;
#  - nudges the sheet/window selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: stamp a destination range with the *exact* cell style already
# used elsewhere in the sheet (copy/paste-special formats only, values
# and formulas in the destination are left untouched). This reuses the
# existing style-table entries instead of synthesizing new ones.
# ---------------------------------------------------------------------
function Copy-Style($srcAddr, [string[]]$destAddrs) {
    $ws.Range($srcAddr).Copy() | Out-Null
    foreach ($d in $destAddrs) {
        $ws.Range($d).PasteSpecial(-4122) | Out-Null
    }
}

# ---------------------------------------------------------------------
# Rows 5, 6, 8, 9, 12, 13: P:T switch from General to "#,##0" (style 7),
# and U switches from bold-General to bold-"#,##0" (style 6).
# ---------------------------------------------------------------------
$dataRows = 5, 6, 8, 9, 12, 13
foreach ($r in $dataRows) {
    $ws.Range("P$r`:T$r").NumberFormat = "#,##0"
    $ws.Range("U$r").NumberFormat = "#,##0"
}

# ---------------------------------------------------------------------
# Rows 7, 10, 11: previously-empty P:T/U cells gain the same formatting
# as the rest of the row (style 7 / style 6), without any values.
# ---------------------------------------------------------------------
$blankRows = 7, 10, 11
foreach ($r in $blankRows) {
    $ws.Range("P$r`:T$r").NumberFormat = "#,##0"
    $ws.Range("U$r").NumberFormat = "#,##0"
}

# ---------------------------------------------------------------------
# Row 13: a new blank A13 cell using the plain "Normal" style (style 0)
# that the rest of column A already uses.
# ---------------------------------------------------------------------
Copy-Style "A2" @("A13")

# ---------------------------------------------------------------------
# Row 14: new placeholder cells mirroring rows 7/10/11 (A/B general,
# F/H/P:T "#,##0", I/U bold "#,##0").
# ---------------------------------------------------------------------
Copy-Style "A2" @("A14", "B14")
$ws.Range("F14").NumberFormat = "#,##0"
$ws.Range("H14").NumberFormat = "#,##0"
$ws.Range("I14").NumberFormat = "#,##0"
$ws.Range("P14:T14").NumberFormat = "#,##0"
$ws.Range("U14").NumberFormat = "#,##0"
# I14 needs the bold "Weight" style (style 6), same as N/I columns elsewhere.
Copy-Style "I5" @("I14")
Copy-Style "N7" @("U14")

# ---------------------------------------------------------------------
# Rows 15 & 16: brand new "alu" function data, columns D:I, matching the
# formulas used throughout the table (F=D+E, H=D+E+G, I=F+(G*5)).
# ---------------------------------------------------------------------
Copy-Style "A2" @("D15", "E15", "G15")
$ws.Range("D15").Value = 41
$ws.Range("E15").Value = 111
$ws.Range("F15").Formula = "=D15+E15"
$ws.Range("G15").Value = 195
$ws.Range("H15").Formula = "=D15+E15+G15"
$ws.Range("I15").Formula = "=F15+(G15*5)"
Copy-Style "F5" @("F15", "H15")
Copy-Style "I5" @("I15")

Copy-Style "A2" @("D16", "E16", "G16")
$ws.Range("D16").Value = 211
$ws.Range("E16").Value = 346
$ws.Range("F16").Formula = "=D16+E16"
$ws.Range("G16").Value = 685
$ws.Range("H16").Formula = "=D16+E16+G16"
$ws.Range("I16").Formula = "=F16+(G16*5)"
Copy-Style "F6" @("F16", "H16")
Copy-Style "I6" @("I16")

# Rows 15/16 P:T / U get the same comma-format upgrade as the other rows.
$ws.Range("P15:T15").NumberFormat = "#,##0"
$ws.Range("U15").NumberFormat = "#,##0"
$ws.Range("P16:T16").NumberFormat = "#,##0"
$ws.Range("U16").NumberFormat = "#,##0"

# ---------------------------------------------------------------------
# Window / sheet view: move the selection down to S25 (the topLeftCell
# scroll position isn't something this headless host persists, but the
# active selection is).
# ---------------------------------------------------------------------
$ws.Range("S25").Select() | Out-Null
